$d = $word.ActiveDocument

$d.Content.Find.Execute("14×46=644", $true, $false, $false, $false, $false, $true, 1, $false, "45×95=4275", 2) | Out-Null
$d.Content.Find.Execute("66×86=5676", $true, $false, $false, $false, $false, $true, 1, $false, "58×50=2900", 2) | Out-Null
$d.Content.Find.Execute("81×11=891", $true, $false, $false, $false, $false, $true, 1, $false, "65×31=2015", 2) | Out-Null
$d.Content.Find.Execute("20×75=1500", $true, $false, $false, $false, $false, $true, 1, $false, "24×77=1848", 2) | Out-Null
$d.Content.Find.Execute("61×80=4880", $true, $false, $false, $false, $false, $true, 1, $false, "98×52=5096", 2) | Out-Null
$d.Content.Find.Execute("96×70=6720", $true, $false, $false, $false, $false, $true, 1, $false, "37×33=1221", 2) | Out-Null
$d.Content.Find.Execute("20×33=660", $true, $false, $false, $false, $false, $true, 1, $false, "67×31=2077", 2) | Out-Null
$d.Content.Find.Execute("48×64=3072", $true, $false, $false, $false, $false, $true, 1, $false, "11×47=517", 2) | Out-Null
$d.Content.Find.Execute("39×94=3666", $true, $false, $false, $false, $false, $true, 1, $false, "75×35=2625", 2) | Out-Null
$d.Content.Find.Execute("87×33=2871", $true, $false, $false, $false, $false, $true, 1, $false, "41×89=3649", 2) | Out-Null
$d.Content.Find.Execute("84×60=5040", $true, $false, $false, $false, $false, $true, 1, $false, "23×11=253", 2) | Out-Null
$d.Content.Find.Execute("57×48=2736", $true, $false, $false, $false, $false, $true, 1, $false, "63×64=4032", 2) | Out-Null
$d.Content.Find.Execute("70×59=4130", $true, $false, $false, $false, $false, $true, 1, $false, "95×32=3040", 2) | Out-Null
$d.Content.Find.Execute("80×66=5280", $true, $false, $false, $false, $false, $true, 1, $false, "97×14=1358", 2) | Out-Null
$d.Content.Find.Execute("19×63=1197", $true, $false, $false, $false, $false, $true, 1, $false, "40×88=3520", 2) | Out-Null
$d.Content.Find.Execute("21×61=1281", $true, $false, $false, $false, $false, $true, 1, $false, "72×11=792", 2) | Out-Null
$d.Content.Find.Execute("48×39=1872", $true, $false, $false, $false, $false, $true, 1, $false, "45×33=1485", 2) | Out-Null
$d.Content.Find.Execute("78×73=5694", $true, $false, $false, $false, $false, $true, 1, $false, "29×16=464", 2) | Out-Null
$d.Content.Find.Execute("10×15=150", $true, $false, $false, $false, $false, $true, 1, $false, "50×23=1150", 2) | Out-Null
$d.Content.Find.Execute("56×28=1568", $true, $false, $false, $false, $false, $true, 1, $false, "60×79=4740", 2) | Out-Null
$d.Content.Find.Execute("57×24=1368", $true, $false, $false, $false, $false, $true, 1, $false, "43×64=2752", 2) | Out-Null
$d.Content.Find.Execute("45×94=4230", $true, $false, $false, $false, $false, $true, 1, $false, "18×22=396", 2) | Out-Null
$d.Content.Find.Execute("84×21=1764", $true, $false, $false, $false, $false, $true, 1, $false, "51×59=3009", 2) | Out-Null
$d.Content.Find.Execute("35×25=875", $true, $false, $false, $false, $false, $true, 1, $false, "93×69=6417", 2) | Out-Null
$d.Content.Find.Execute("20×15=300", $true, $false, $false, $false, $false, $true, 1, $false, "100×70=7000", 2) | Out-Null
$d.Content.Find.Execute("91×52=4732", $true, $false, $false, $false, $false, $true, 1, $false, "33×95=3135", 2) | Out-Null
$d.Content.Find.Execute("34×32=1088", $true, $false, $false, $false, $false, $true, 1, $false, "11×21=231", 2) | Out-Null
$d.Content.Find.Execute("98×13=1274", $true, $false, $false, $false, $false, $true, 1, $false, "90×99=8910", 2) | Out-Null
$d.Content.Find.Execute("90×33=2970", $true, $false, $false, $false, $false, $true, 1, $false, "69×54=3726", 2) | Out-Null
$d.Content.Find.Execute("87×38=3306", $true, $false, $false, $false, $false, $true, 1, $false, "38×60=2280", 2) | Out-Null
$d.Content.Find.Execute("47×60=2820", $true, $false, $false, $false, $false, $true, 1, $false, "64×27=1728", 2) | Out-Null
$d.Content.Find.Execute("76×96=7296", $true, $false, $false, $false, $false, $true, 1, $false, "75×23=1725", 2) | Out-Null
$d.Content.Find.Execute("62×15=930", $true, $false, $false, $false, $false, $true, 1, $false, "54×56=3024", 2) | Out-Null
$d.Content.Find.Execute("77×63=4851", $true, $false, $false, $false, $false, $true, 1, $false, "85×92=7820", 2) | Out-Null
$d.Content.Find.Execute("96×86=8256", $true, $false, $false, $false, $false, $true, 1, $false, "83×55=4565", 2) | Out-Null
$d.Content.Find.Execute("45×34=1530", $true, $false, $false, $false, $false, $true, 1, $false, "63×84=5292", 2) | Out-Null
$d.Content.Find.Execute("69×53=3657", $true, $false, $false, $false, $false, $true, 1, $false, "81×33=2673", 2) | Out-Null
$d.Content.Find.Execute("48×52=2496", $true, $false, $false, $false, $false, $true, 1, $false, "78×85=6630", 2) | Out-Null
$d.Content.Find.Execute("95×85=8075", $true, $false, $false, $false, $false, $true, 1, $false, "43×96=4128", 2) | Out-Null
$d.Content.Find.Execute("27×35=945", $true, $false, $false, $false, $false, $true, 1, $false, "86×56=4816", 2) | Out-Null
$d.Content.Find.Execute("26×47=1222", $true, $false, $false, $false, $false, $true, 1, $false, "25×62=1550", 2) | Out-Null
$d.Content.Find.Execute("81×89=7209", $true, $false, $false, $false, $false, $true, 1, $false, "58×82=4756", 2) | Out-Null
$d.Content.Find.Execute("95×14=1330", $true, $false, $false, $false, $false, $true, 1, $false, "76×17=1292", 2) | Out-Null
$d.Content.Find.Execute("86×79=6794", $true, $false, $false, $false, $false, $true, 1, $false, "81×72=5832", 2) | Out-Null
$d.Content.Find.Execute("47×85=3995", $true, $false, $false, $false, $false, $true, 1, $false, "22×54=1188", 2) | Out-Null
$d.Content.Find.Execute("91×24=2184", $true, $false, $false, $false, $false, $true, 1, $false, "81×15=1215", 2) | Out-Null
$d.Content.Find.Execute("62×24=1488", $true, $false, $false, $false, $false, $true, 1, $false, "41×46=1886", 2) | Out-Null
$d.Content.Find.Execute("58×40=2320", $true, $false, $false, $false, $false, $true, 1, $false, "35×40=1400", 2) | Out-Null
$d.Content.Find.Execute("55×10=550", $true, $false, $false, $false, $false, $true, 1, $false, "62×14=868", 2) | Out-Null
$d.Content.Find.Execute("45×90=4050", $true, $false, $false, $false, $false, $true, 1, $false, "97×80=7760", 2) | Out-Null
$d.Content.Find.Execute("32×35=1120", $true, $false, $false, $false, $false, $true, 1, $false, "96×18=1728", 2) | Out-Null
$d.Content.Find.Execute("57×16=912", $true, $false, $false, $false, $false, $true, 1, $false, "27×30=810", 2) | Out-Null
$d.Content.Find.Execute("89×19=1691", $true, $false, $false, $false, $false, $true, 1, $false, "57×71=4047", 2) | Out-Null
$d.Content.Find.Execute("41×52=2132", $true, $false, $false, $false, $false, $true, 1, $false, "43×17=731", 2) | Out-Null
$d.Content.Find.Execute("12×13=156", $true, $false, $false, $false, $false, $true, 1, $false, "70×82=5740", 2) | Out-Null
$d.Content.Find.Execute("61×50=3050", $true, $false, $false, $false, $false, $true, 1, $false, "94×98=9212", 2) | Out-Null
$d.Content.Find.Execute("95×33=3135", $true, $false, $false, $false, $false, $true, 1, $false, "22×19=418", 2) | Out-Null
$d.Content.Find.Execute("53×75=3975", $true, $false, $false, $false, $false, $true, 1, $false, "84×68=5712", 2) | Out-Null
$d.Content.Find.Execute("47×23=1081", $true, $false, $false, $false, $false, $true, 1, $false, "87×58=5046", 2) | Out-Null
$d.Content.Find.Execute("10×36=360", $true, $false, $false, $false, $false, $true, 1, $false, "19×97=1843", 2) | Out-Null
$d.Content.Find.Execute("15×72=1080", $true, $false, $false, $false, $false, $true, 1, $false, "19×88=1672", 2) | Out-Null
$d.Content.Find.Execute("37×26=962", $true, $false, $false, $false, $false, $true, 1, $false, "16×79=1264", 2) | Out-Null
$d.Content.Find.Execute("87×43=3741", $true, $false, $false, $false, $false, $true, 1, $false, "100×67=6700", 2) | Out-Null
$d.Content.Find.Execute("69×20=1380", $true, $false, $false, $false, $false, $true, 1, $false, "61×30=1830", 2) | Out-Null
$d.Content.Find.Execute("79×84=6636", $true, $false, $false, $false, $false, $true, 1, $false, "100×64=6400", 2) | Out-Null
$d.Content.Find.Execute("28×63=1764", $true, $false, $false, $false, $false, $true, 1, $false, "22×92=2024", 2) | Out-Null
$d.Content.Find.Execute("72×80=5760", $true, $false, $false, $false, $false, $true, 1, $false, "14×30=420", 2) | Out-Null
$d.Content.Find.Execute("86×84=7224", $true, $false, $false, $false, $false, $true, 1, $false, "15×87=1305", 2) | Out-Null
$d.Content.Find.Execute("60×67=4020", $true, $false, $false, $false, $false, $true, 1, $false, "90×47=4230", 2) | Out-Null
$d.Content.Find.Execute("17×21=357", $true, $false, $false, $false, $false, $true, 1, $false, "51×33=1683", 2) | Out-Null
$d.Content.Find.Execute("46×93=4278", $true, $false, $false, $false, $false, $true, 1, $false, "15×56=840", 2) | Out-Null
$d.Content.Find.Execute("36×65=2340", $true, $false, $false, $false, $false, $true, 1, $false, "75×25=1875", 2) | Out-Null
$d.Content.Find.Execute("100×86=8600", $true, $false, $false, $false, $false, $true, 1, $false, "12×66=792", 2) | Out-Null
$d.Content.Find.Execute("52×52=2704", $true, $false, $false, $false, $false, $true, 1, $false, "98×67=6566", 2) | Out-Null
$d.Content.Find.Execute("68×40=2720", $true, $false, $false, $false, $false, $true, 1, $false, "28×96=2688", 2) | Out-Null
$d.Content.Find.Execute("80×38=3040", $true, $false, $false, $false, $false, $true, 1, $false, "81×39=3159", 2) | Out-Null
$d.Content.Find.Execute("30×28=840", $true, $false, $false, $false, $false, $true, 1, $false, "51×100=5100", 2) | Out-Null
$d.Content.Find.Execute("63×85=5355", $true, $false, $false, $false, $false, $true, 1, $false, "46×66=3036", 2) | Out-Null
$d.Content.Find.Execute("87×31=2697", $true, $false, $false, $false, $false, $true, 1, $false, "67×85=5695", 2) | Out-Null
$d.Content.Find.Execute("31×10=310", $true, $false, $false, $false, $false, $true, 1, $false, "52×87=4524", 2) | Out-Null
$d.Content.Find.Execute("85×52=4420", $true, $false, $false, $false, $false, $true, 1, $false, "47×83=3901", 2) | Out-Null
$d.Content.Find.Execute("88×37=3256", $true, $false, $false, $false, $false, $true, 1, $false, "86×88=7568", 2) | Out-Null
$d.Content.Find.Execute("28×20=560", $true, $false, $false, $false, $false, $true, 1, $false, "40×35=1400", 2) | Out-Null
$d.Content.Find.Execute("32×48=1536", $true, $false, $false, $false, $false, $true, 1, $false, "65×45=2925", 2) | Out-Null
$d.Content.Find.Execute("56×54=3024", $true, $false, $false, $false, $false, $true, 1, $false, "28×99=2772", 2) | Out-Null
$d.Content.Find.Execute("29×25=725", $true, $false, $false, $false, $false, $true, 1, $false, "25×44=1100", 2) | Out-Null
$d.Content.Find.Execute("17×70=1190", $true, $false, $false, $false, $false, $true, 1, $false, "86×81=6966", 2) | Out-Null
$d.Content.Find.Execute("25×40=1000", $true, $false, $false, $false, $false, $true, 1, $false, "71×55=3905", 2) | Out-Null
$d.Content.Find.Execute("65×79=5135", $true, $false, $false, $false, $false, $true, 1, $false, "24×60=1440", 2) | Out-Null
$d.Content.Find.Execute("50×33=1650", $true, $false, $false, $false, $false, $true, 1, $false, "33×60=1980", 2) | Out-Null
$d.Content.Find.Execute("67×48=3216", $true, $false, $false, $false, $false, $true, 1, $false, "13×18=234", 2) | Out-Null
$d.Content.Find.Execute("39×28=1092", $true, $false, $false, $false, $false, $true, 1, $false, "16×54=864", 2) | Out-Null
$d.Content.Find.Execute("62×71=4402", $true, $false, $false, $false, $false, $true, 1, $false, "87×27=2349", 2) | Out-Null
$d.Content.Find.Execute("37×85=3145", $true, $false, $false, $false, $false, $true, 1, $false, "91×26=2366", 2) | Out-Null
$d.Content.Find.Execute("71×13=923", $true, $false, $false, $false, $false, $true, 1, $false, "10×31=310", 2) | Out-Null
$d.Content.Find.Execute("80×97=7760", $true, $false, $false, $false, $false, $true, 1, $false, "32×72=2304", 2) | Out-Null
$d.Content.Find.Execute("95×13=1235", $true, $false, $false, $false, $false, $true, 1, $false, "49×36=1764", 2) | Out-Null
$d.Content.Find.Execute("33×67=2211", $true, $false, $false, $false, $false, $true, 1, $false, "12×94=1128", 2) | Out-Null
$d.Content.Find.Execute("80×75=6000", $true, $false, $false, $false, $false, $true, 1, $false, "28×45=1260", 2) | Out-Null
$d.Content.Find.Execute("24×34=816", $true, $false, $false, $false, $false, $true, 1, $false, "60×20=1200", 2) | Out-Null
